$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 98, pushing the existing rows 98:119 down to 99:120
$ws.Rows("98:98").Insert()

# Populate the newly inserted row 98 with the new record's data.
# (Columns A,B,C,E,F,G,H,I,J,K,L,Q,R,T are constant across every Mango /
# Macroferia Regional de Talca row in this sheet.)
$ws.Range("A98").Value = 5
$ws.Range("B98").Value = "Macroferia Regional de Talca"
$ws.Range("C98").Value = "Maule"
$ws.Range("D98").Value = 44637
$ws.Range("E98").Value = 7
$ws.Range("F98").Value = "Fruta"
$ws.Range("G98").Value = 100108
$ws.Range("H98").Value = "Tropicales y subtropicales"
$ws.Range("I98").Value = 100108002
$ws.Range("J98").Value = "Mango"
$ws.Range("K98").Value = "Sin especificar"
$ws.Range("L98").Value = "Primera"
$ws.Range("M98").Value = 200
$ws.Range("N98").Value = 7000
$ws.Range("O98").Value = 7000
$ws.Range("P98").Value = 7000
$ws.Range("Q98").Value = '$/bandeja 4 kilos'
$ws.Range("R98").Value = "Perú"
$ws.Range("S98").Value = 1750
$ws.Range("T98").Value = 4
